# Raw and Clean Data from SSA for July 30th
# Appends the 2020-07-30 (serial 44042) data row/column to each tracking
# sheet in the daily bitacora workbook, and clears the stale preview
# column on control_obs_mpio.

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# out_vars: brand-new row 61 (summary of the day's raw/clean counts)
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("out_vars")
$ws.Range("A60").Copy()
$ws.Range("A61").PasteSpecial(-4122)
$ws.Range("B60").Copy()
$ws.Range("B61:J61").PasteSpecial(-4122)

$ws.Range("A61").Value = 44042
$ws.Range("B61").Value = 416179
$ws.Range("C61").Value = 461775
$ws.Range("D61").Value = 90582
$ws.Range("E61").Value = 46000
$ws.Range("F61").Value = 27.380045605376534
$ws.Range("G61").Value = 113950
$ws.Range("H61").Value = 9188
$ws.Range("I61").Value = 10909
$ws.Range("J61").Value = 968536

# ----------------------------------------------------------------------
# dates_dx: the placeholder row 61 gets real counts
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("dates_dx")
$wb.Worksheets.Item("out_vars").Range("A60").Copy()
$ws.Range("A61").PasteSpecial(-4122)

$ws.Range("A61").Value = 44042
$ws.Range("B61").Value = 0
$ws.Range("C61").Value = 1
$ws.Range("D61").Value = 0
$ws.Range("E61").Value = 0
$ws.Range("F61").Value = 1
$ws.Range("G61").Value = 0
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 4

# ----------------------------------------------------------------------
# dates_sx: row 61 previously only had a blank A61; fill the whole row
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("dates_sx")
$wb.Worksheets.Item("out_vars").Range("A60").Copy()
$ws.Range("A61").PasteSpecial(-4122)
$ws.Range("C58").Copy()
$ws.Range("B61:N61").PasteSpecial(-4122)

$ws.Range("A61").Value = 44042
$ws.Range("B61").Value = 0
$ws.Range("C61").Value = 1
$ws.Range("D61").Value = 0
$ws.Range("E61").Value = 0
$ws.Range("F61").Value = 0
$ws.Range("G61").Value = 0
$ws.Range("H61").Value = 1
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = 0
$ws.Range("N61").Value = 0

# ----------------------------------------------------------------------
# dates_deaths: placeholder row 61 gets real counts
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("dates_deaths")
$wb.Worksheets.Item("out_vars").Range("A60").Copy()
$ws.Range("A61").PasteSpecial(-4122)

$ws.Range("A61").Value = 44042
$ws.Range("B61").Value = 0
$ws.Range("C61").Value = 0
$ws.Range("D61").Value = 0
$ws.Range("E61").Value = 0
$ws.Range("F61").Value = 2
$ws.Range("G61").Value = 1
$ws.Range("H61").Value = 1
$ws.Range("I61").Value = 1
$ws.Range("J61").Value = 2

# ----------------------------------------------------------------------
# control_obs: new column BI (the 7/30 snapshot column)
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("control_obs")
$ws.Range("BH1").Copy()
$ws.Range("BI1").PasteSpecial(-4122)
$ws.Range("BH18").Copy()
$ws.Range("BI18").PasteSpecial(-4122)
$ws.Range("C2").Copy()
$ws.Range("BI20").PasteSpecial(-4122)

$ws.Range("BI1").Value = 44042
$ws.Range("BI2").Value = 4834
$ws.Range("BI3").Value = 4633
$ws.Range("BI4").Value = 4633
$ws.Range("BI5").Value = 4633
$ws.Range("BI6").Value = 4633
$ws.Range("BI7").Value = 3943
$ws.Range("BI8").Value = 6574
$ws.Range("BI10").Value = 200
$ws.Range("BI11").Value = 200
$ws.Range("BI12").Value = 200
$ws.Range("BI13").Value = 200
$ws.Range("BI14").Value = 200
$ws.Range("BI15").Value = 135
$ws.Range("BI16").Value = 212
$ws.Range("BI18").Value = 1129
$ws.Range("BI20").Formula = "=SUM(BI2:BI18)"

# ----------------------------------------------------------------------
# control_obs_mpio: drop the stale preview column G (superseded by the
# new BI column on control_obs)
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("control_obs_mpio")
$ws.Range("G1:G10").ClearContents()

# ----------------------------------------------------------------------
# view / selection bookkeeping to mirror the author's final UI state
# ----------------------------------------------------------------------
$wb.Worksheets.Item("out_vars").Range("C70").Select()
$wb.Worksheets.Item("dates_dx").Range("L61").Select()
$wb.Worksheets.Item("dates_sx").Range("O61").Select()
$wb.Worksheets.Item("dates_deaths").Range("K61").Select()
$wb.Worksheets.Item("control_obs").Range("BH29").Select()
$wb.Worksheets.Item("control_obs_mpio").Range("H21").Select()
$wb.Worksheets.Item("anomalias").Activate()
$wb.Worksheets.Item("out_vars").Activate()
